# Remove the duplicated slide at position 6 ("VACCINES : APPROVALS")
$p = $ppt.ActivePresentation
$p.Slides.Item(6).Delete()
